$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J16").Value = "88-90 MPH"
$ws.Range("J17").Value = "CB,FB,CH"

$ws.Range("J25").Value = "88-90 MPH"
$ws.Range("J26").Value = "CB,FB,CH"

$ws.Range("J28").Value = 5
$ws.Range("M28").Value = "82.59 MPH"
$ws.Range("J29").Value = 2
$ws.Range("M30").Value = "-15.88°"

$ws.Range("J32").Value = "Herbst"
$ws.Range("M32").Value = "Ground Ball"
$ws.Range("M33").Value = "Double"
$ws.Range("J34").Value = "83-85 MPH"
$ws.Range("J35").Value = "SL,CB,FB,CH"

$ws.Range("J37").Value = 3
$ws.Range("M37").Value = "83.5 MPH"
$ws.Range("J38").Value = 1
$ws.Range("M39").Value = "32.75°"

$ws.Range("J41").Value = "Roblez"
$ws.Range("M41").Value = "Fly Ball"
$ws.Range("M42").Value = "Out"
$ws.Range("J43").Value = "88-90 MPH"
$ws.Range("J44").Value = "CB,FB,CH"
